$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for this market/product. It belongs
# right after the existing row 247, so insert a fresh row at 248 — this
# shifts the old rows 248:270 down to 249:271 (carrying their formatting
# along), growing the used range from A1:T270 to A1:T271.
$ws.Rows("248:248").Insert()

# Populate the newly inserted row 248 with the new weekly observation.
# Most fields repeat the constant market/product metadata used throughout
# the sheet; only the date and the price columns (N, O, P, S) are new.
$ws.Range("A248").Value = 1
$ws.Range("B248").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C248").Value = "Arica y Parinacota"
$ws.Range("D248").Value = 44769
$ws.Range("E248").Value = 15
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100108
$ws.Range("H248").Value = "Tropicales y subtropicales"
$ws.Range("I248").Value = 100108006
$ws.Range("J248").Value = "Plátano"
$ws.Range("K248").Value = "Sin especificar"
$ws.Range("L248").Value = "Pintón"
$ws.Range("M248").Value = 120
$ws.Range("N248").Value = 27000
$ws.Range("O248").Value = 28000
$ws.Range("P248").Value = 27500
$ws.Range("Q248").Value = '$/caja 20 kilos'
$ws.Range("R248").Value = "Ecuador"
$ws.Range("S248").Value = 1375
$ws.Range("T248").Value = 20
